$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 695.3077
$ws.Range("I6").Value = 130.81818
$ws.Range("J6").Value = 3800
$ws.Range("K6").Value = 392.4545400000001
$ws.Range("L6").Value = 11400
$ws.Range("M6").Value = -280.4545400000001
$ws.Range("N6").Value = -11624

$ws.Range("H13").Value = 25002.75
$ws.Range("J13").Value = 25002.75
$ws.Range("L13").Value = 25002.75
$ws.Range("N13").Value = -25340.75

$ws.Range("H18").Value = 320.5625
$ws.Range("I18").Value = 207.27272
$ws.Range("J18").Value = 569.8
$ws.Range("K18").Value = 207.27272
$ws.Range("L18").Value = 569.8
$ws.Range("M18").Value = 76.72728000000001
$ws.Range("N18").Value = -1137.8

$ws.Range("H19").Value = 2025629.1
$ws.Range("I19").Value = 3760390.8
$ws.Range("K19").Value = 3760390.8
$ws.Range("M19").Value = -3760215.8

$ws.Range("H41").Value = 950
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 1071.4286
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 1071.4286
$ws.Range("M41").Value = 340
$ws.Range("N41").Value = -1951.4286

$ws.Range("H42").Value = 542
$ws.Range("I42").Value = 468.625
$ws.Range("J42").Value = 595.36365
$ws.Range("K42").Value = 1405.875
$ws.Range("L42").Value = 1786.09095
$ws.Range("M42").Value = -1175.875
$ws.Range("N42").Value = -2246.09095

$ws.Range("H43").Value = 1217.5
$ws.Range("I43").Value = 945
$ws.Range("J43").Value = 1278.0555
$ws.Range("K43").Value = 945
$ws.Range("L43").Value = 1278.0555
$ws.Range("M43").Value = -876
$ws.Range("N43").Value = -1416.0555

$ws.Range("H51").Value = 7033.3335
$ws.Range("I51").Value = 3000
$ws.Range("K51").Value = 3000
$ws.Range("M51").Value = -2516

$ws.Range("H53").Value = 740.4666999999999
$ws.Range("J53").Value = 725.5833
$ws.Range("L53").Value = 725.5833
$ws.Range("N53").Value = -1999.5833

$ws.Range("H55").Value = 300
$ws.Range("I55").Value = 300
$ws.Range("K55").Value = 300
$ws.Range("M55").Value = -86

$ws.Range("H98").Value = 3608.4167
$ws.Range("I98").Value = 1872.7222
$ws.Range("J98").Value = 8815.5
$ws.Range("K98").Value = 1872.7222
$ws.Range("L98").Value = 8815.5
$ws.Range("M98").Value = -374.7221999999999
$ws.Range("N98").Value = -11811.5

$ws.Range("H122").Value = 3608.4167
$ws.Range("I122").Value = 1872.7222
$ws.Range("J122").Value = 8815.5
$ws.Range("K122").Value = 5618.1666
$ws.Range("L122").Value = 26446.5
$ws.Range("M122").Value = -3168.1666
$ws.Range("N122").Value = -31346.5

$ws.Range("H129").Value = 1036.9354
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 1059.8334
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 3179.5002
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -13179.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 38712
$ws.Range("J7").Value = 38712
$ws.Range("L7").Value = 38712
$ws.Range("N7").Value = -38940

$ws.Range("H32").Value = 6125.7964
$ws.Range("I32").Value = 3757.1714
$ws.Range("K32").Value = 3757.1714
$ws.Range("M32").Value = -3470.1714

$ws.Range("H122").Value = 2874.1177
$ws.Range("I122").Value = 1800
$ws.Range("K122").Value = 5400
$ws.Range("M122").Value = -2950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 199
$ws.Range("I22").Value = 199
$ws.Range("K22").Value = 199
$ws.Range("M22").Value = -26

$ws.Range("H94").Value = 1249.1333
$ws.Range("I94").Value = 906.2083
$ws.Range("K94").Value = 906.2083
$ws.Range("M94").Value = -455.2083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11768559
$ws.Range("I99").Value = 20002230
$ws.Range("J99").Value = 6172
$ws.Range("K99").Value = 20002230
$ws.Range("L99").Value = 6172
$ws.Range("M99").Value = -20000732
$ws.Range("N99").Value = -9168

$ws.Range("H105").Value = 1707.25
$ws.Range("I105").Value = 1555.5333
$ws.Range("J105").Value = 1960.1111
$ws.Range("K105").Value = 1555.5333
$ws.Range("L105").Value = 1960.1111
$ws.Range("M105").Value = 191.4666999999999
$ws.Range("N105").Value = -5454.1111

$ws.Range("H122").Value = 2471.875
$ws.Range("I122").Value = 1426
$ws.Range("J122").Value = 4215
$ws.Range("K122").Value = 4278
$ws.Range("L122").Value = 12645
$ws.Range("M122").Value = -1828
$ws.Range("N122").Value = -17545

$ws.Range("H126").Value = 11768559
$ws.Range("I126").Value = 20002230
$ws.Range("J126").Value = 6172
$ws.Range("K126").Value = 60006690
$ws.Range("L126").Value = 18516
$ws.Range("M126").Value = -60004220
$ws.Range("N126").Value = -23456

$ws.Range("H132").Value = 3819.1724
$ws.Range("I132").Value = 4346
$ws.Range("J132").Value = 3391.125
$ws.Range("K132").Value = 13038
$ws.Range("L132").Value = 10173.375
$ws.Range("M132").Value = -10508
$ws.Range("N132").Value = -15233.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 144.55556
$ws.Range("J38").Value = 178.5
$ws.Range("L38").Value = 535.5
$ws.Range("N38").Value = -1229.5

$ws.Range("H116").Value = 3666.6667
$ws.Range("I116").Value = 1000
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -21884

$ws.Range("H131").Value = 8065543
$ws.Range("J131").Value = 911.63794
$ws.Range("L131").Value = 2734.91382
$ws.Range("N131").Value = -12814.91382

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6947239
$ws.Range("I80").Value = 10872266
$ws.Range("J80").Value = 2961.5386
$ws.Range("K80").Value = 10872266
$ws.Range("L80").Value = 2961.5386
$ws.Range("M80").Value = -10871268
$ws.Range("N80").Value = -4957.5386

$ws.Range("H83").Value = 6947239
$ws.Range("I83").Value = 10872266
$ws.Range("J83").Value = 2961.5386
$ws.Range("K83").Value = 54361330
$ws.Range("L83").Value = 14807.693
$ws.Range("M83").Value = -54356338
$ws.Range("N83").Value = -24791.693

$ws.Range("H120").Value = 27654.25
$ws.Range("J120").Value = 27654.25
$ws.Range("L120").Value = 27654.25
$ws.Range("N120").Value = -37330.25

$ws.Range("H134").Value = 42907.363
$ws.Range("J134").Value = 42907.363
$ws.Range("L134").Value = 128722.089
$ws.Range("N134").Value = -133792.089

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4495.9565
$ws.Range("I7").Value = 2684.3333
$ws.Range("K7").Value = 2684.3333
$ws.Range("M7").Value = -2572.3333

$ws.Range("H22").Value = 2880.875
$ws.Range("I22").Value = 2165.3333
$ws.Range("K22").Value = 2165.3333
$ws.Range("M22").Value = -1870.3333

$ws.Range("H27").Value = 2880.875
$ws.Range("I27").Value = 2165.3333
$ws.Range("K27").Value = 2165.3333
$ws.Range("M27").Value = -2058.3333

$ws.Range("H122").Value = 6535.2856
$ws.Range("I122").Value = 3643.4285
$ws.Range("K122").Value = 10930.2855
$ws.Range("M122").Value = -8480.2855

$ws.Range("H126").Value = 4495.9565
$ws.Range("I126").Value = 2684.3333
$ws.Range("K126").Value = 8052.999899999999
$ws.Range("M126").Value = -5582.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3196.5356
$ws.Range("I122").Value = 1715.9474
$ws.Range("J122").Value = 6322.222
$ws.Range("K122").Value = 5147.8422
$ws.Range("L122").Value = 18966.666
$ws.Range("M122").Value = -2697.8422
$ws.Range("N122").Value = -23866.666

$ws.Range("H136").Value = 1699.8214
$ws.Range("I136").Value = 1182.3334
$ws.Range("J136").Value = 2296.923
$ws.Range("K136").Value = 3547.0002
$ws.Range("L136").Value = 6890.768999999999
$ws.Range("M136").Value = -997.0001999999999
$ws.Range("N136").Value = -11990.769
